# Value of a Statistical Life.xlsx - "Uploading newest US files"
#
# The "About" sheet's Notes section previously explained the VSL variable's
# relationship to the SCoHIbP (Social Cost of Health Impacts by Pollutant)
# variable across four lines (rows 10-13) followed by a blank spacer row
# (row 14). That note is replaced with a shorter, two-line explanation
# referencing HOIpTP instead, and the spacer row collapses to a single gap.
# Everything below shifts up to follow.

$wb = $excel.ActiveWorkbook

$about = $wb.Worksheets.Item("About")
$voasl = $wb.Worksheets.Item("VoaSL")

# Remove the old 4-line SCoHIbP note plus the blank spacer line beneath it
# (old rows 10-14), then reopen 3 fresh rows for the new, shorter note plus
# its own single blank spacer line.
$about.Range("A10:A14").EntireRow.Delete()
$about.Range("A10:A12").EntireRow.Insert()

$about.Range("A10").Value = "This variable is used to convert estimated avoided premature mortalities"
$about.Range("A11").Value = "(calculated using data from HOIpTP) into a dollar amount."

# The inserted rows pick up formatting from the row above; the replacement
# note text isn't bold like the "Notes" heading, so plain font is restored.
$about.Range("A10:A12").Font.Bold = $false

# Row 12 is a blank spacer (no styling) - fully clear it so it matches the
# other blank spacer rows on this sheet (e.g. no leftover cell record).
$about.Range("A12").Clear()

# The VoaSL sheet's conversion formula references the CPI adjustment factor
# on the About sheet; that factor now lives two rows higher (About!A15
# instead of About!A17) after the note shrank.
$voasl.Range("B2").Formula = "=7.4*10^6*About!A15"

# Restore the active selection on the About sheet.
$about.Activate()
$about.Range("I18").Select() | Out-Null
